$d = $word.ActiveDocument

# The diff swaps the East-Asian fallback font from "DejaVu Sans" to
# "Tahoma" for the document defaults / Normal / Heading styles, and it
# adds an explicit complex-script ("cs") font of "DejaVu Sans" to the
# List, Caption and Index styles (which previously had no rFonts entry
# at all, relying on inheritance).

# Normal style: eastAsia DejaVu Sans -> Tahoma
$normal = $d.Styles("Normal")
$normal.Font.NameFarEast = "Tahoma"

# Heading style: eastAsia DejaVu Sans -> Tahoma
$heading = $d.Styles("Heading")
$heading.Font.NameFarEast = "Tahoma"

# List style: add <w:rFonts w:cs="DejaVu Sans"/>
$list = $d.Styles("List")
$list.Font.NameBi = "DejaVu Sans"

# Caption style: add <w:rFonts w:cs="DejaVu Sans"/> (kept before i/iCs)
$caption = $d.Styles("Caption")
$caption.Font.NameBi = "DejaVu Sans"

# Index style: add <w:rFonts w:cs="DejaVu Sans"/>
$index = $d.Styles("Index")
$index.Font.NameBi = "DejaVu Sans"
